$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B85").Value = "JS101: Programming Foundations with JavaScript"
$ws.Range("C85").Value = 1
$ws.Range("D85").Value = "Finish 4 small problems"
